$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Orange"
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = "Tomato"
$ws.Range("B6").Value = 3

$ws.Range("A6").Select()
